$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cesium renderer status update: commit "Add PixelFormat, start implementing renderer texture"
#
# For each of the rows below, the old "Status" column (B) is removed
# entirely (it duplicated the "Swift port" info) and the "Swift port"
# column (C) is updated to reflect the real per-file swift-port state.
# Row 208 (RenderState.js) also gets a new note in column E.

# Module.js                 -> new Swift port (C) value
# Row 92  Matrix4.js         -> Partial
$ws.Range("C92").Value = "Partial"
$ws.Range("B92").ClearContents()

# Row 198 Context.js         -> Complete
$ws.Range("C198").Value = "Complete"
$ws.Range("B198").ClearContents()

# Row 199 CubeMap.js         -> Partial
$ws.Range("C199").Value = "Partial"
$ws.Range("B199").ClearContents()

# Row 200 CubeMapFace.js     -> Partial
$ws.Range("C200").Value = "Partial"
$ws.Range("B200").ClearContents()

# Row 202 DrawCommand.js     -> Complete
$ws.Range("C202").Value = "Complete"
$ws.Range("B202").ClearContents()

# Row 203 Framebuffer.js     -> Partial
$ws.Range("C203").Value = "Partial"
$ws.Range("B203").ClearContents()

# Row 205 PickFramebuffer.js -> Partial
$ws.Range("C205").Value = "Partial"
$ws.Range("B205").ClearContents()

# Row 208 RenderState.js     -> Partial, note changes to "Needs function array"
$ws.Range("C208").Value = "Partial"
$ws.Range("B208").ClearContents()
$ws.Range("E208").Value = "Needs function array"

# Row 212 ShaderCache.js     -> Partial
$ws.Range("C212").Value = "Partial"
$ws.Range("B212").ClearContents()

# Row 213 ShaderProgram.js   -> Partial
$ws.Range("C213").Value = "Partial"
$ws.Range("B213").ClearContents()

# Row 214 Texture.js         -> Partial
$ws.Range("C214").Value = "Partial"
$ws.Range("B214").ClearContents()

# Row 218 UniformState.js    -> Partial
$ws.Range("C218").Value = "Partial"
$ws.Range("B218").ClearContents()

# Row 219 VertexArray.js     -> Partial
$ws.Range("C219").Value = "Partial"
$ws.Range("B219").ClearContents()

# Reflect the author's on-screen position/selection at save time: scrolled
# down to row 80, with B92 selected.
$excel.ActiveWindow.ScrollRow = 80
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B92").Select()
